$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Apparel and Clothing"
$ws.Range("B7").Value = 13
$ws.Range("C7").Value = 681

$ws.Range("A7").Select()
